$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (AGN names re-shuffled / epochs corrected) ---
$ws.Range("A5").Value = "090159.65+333551.1"
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = "123101.16+555023.5"
$ws.Range("C6").Value = 4

$ws.Range("C7").Value = 9

$ws.Range("A8").Value = "142030.56+513909.1"
$ws.Range("C8").Value = 6

$ws.Range("A9").Value = "142320.98+540509.3"
$ws.Range("C9").Value = 5

$ws.Range("A10").Value = "163618.68+331916.9"
$ws.Range("C10").Value = 9

$ws.Range("A11").Value = "123632.75+552109.3"
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = "New"

$ws.Range("A12").Value = "123632.75+552109.3"
$ws.Range("C12").Value = 9

# --- Add new rows 13-23 ---
$ws.Range("A13").Value = "122321.78+311410.6"
$ws.Range("B13").Value = "W1"
$ws.Range("C13").Value = 3

$ws.Range("A14").Value = "151354.51+184804.3"
$ws.Range("B14").Value = "W1"
$ws.Range("C14").Value = 8

$ws.Range("A15").Value = "124939.87+581115.5"
$ws.Range("B15").Value = "W1"
$ws.Range("C15").Value = 14

$ws.Range("A16").Value = "155021.33+432712.6"
$ws.Range("B16").Value = "W1"
$ws.Range("C16").Value = 12

$ws.Range("A17").Value = "155021.33+432712.6"
$ws.Range("B17").Value = "W1"
$ws.Range("C17").Value = 13

$ws.Range("A18").Value = "111353.73+515725.8"
$ws.Range("B18").Value = "W1"
$ws.Range("C18").Value = 4

$ws.Range("A19").Value = "125731.87+272313.3"
$ws.Range("B19").Value = "W1"
$ws.Range("C19").Value = 7

$ws.Range("A20").Value = "125731.87+272313.3"
$ws.Range("B20").Value = "W1"
$ws.Range("C20").Value = 13

$ws.Range("A21").Value = "125731.87+272313.3"
$ws.Range("B21").Value = "W1"
$ws.Range("C21").Value = 14

$ws.Range("A22").Value = "162659.42+424450.0"
$ws.Range("B22").Value = "W1"
$ws.Range("C22").Value = 14

$ws.Range("A23").Value = "162659.42+424450.0"
$ws.Range("B23").Value = "W1"
$ws.Range("C23").Value = 14

# --- Update selection to match final saved state ---
$ws.Range("F7").Select()
